# Aceite de Canola (row 9) currently has its Fat content (column D) entered
# as the text "100.0%" (a quote-prefixed / text-formatted percentage). The
# edit replaces that with the real numeric value 1 (i.e. 100%), formatted
# the same way as the other percentage cells already in column D (e.g. D4),
# so it reads as "100.0%" via a genuine percentage number format instead of
# literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing percentage formatting (font/fill/border/alignment/
# number format) already applied to the other Fat-content cells in column D
# by copying D4's format onto D9, then overwrite the value with the real
# number 1 (100%).
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D9").Value = 1

# Match the resulting selection state (the user ended up with D9 selected).
$ws.Range("D9").Select() | Out-Null
